$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reviewer comments: phenol measurements unit revised. The "D" column
# (Phenols_level, mg/100g OC) stays as-is; a new "E" column reports the
# same quantity on a per-100g-OC basis scaled by 10 and is relabeled
# "Phenols_level" (replacing the old SD column).
$ws.Range("D1").Value = "Phenols_level_mg_100gOC"
$ws.Range("E1").Value = "Phenols_level"

for ($r = 2; $r -le 21; $r++) {
    $d = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 5).Value = $d * 10
}

# Resize the header columns to fit the new, longer header text.
$ws.Columns("D").ColumnWidth = 24.5
$ws.Columns("E").ColumnWidth = 12.833333333333334

# Selection narrows from E2:E5 down to just E2.
$ws.Range("E2").Select() | Out-Null
